# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# 展览 sheet:    F4 2062 -> 2067, F5 1650 -> 1653, F8 673 -> 677
# 全部类型 sheet: F4 2062 -> 2067, F5 1650 -> 1653, F9 673 -> 677

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2067
$ws1.Range("F5").Value = 1653
$ws1.Range("F8").Value = 677

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2067
$ws4.Range("F5").Value = 1653
$ws4.Range("F9").Value = 677
